$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsAttr = $wb.Worksheets.Item("Attribute description")

# ---------------------------------------------------------------------------
# 1. isolation_source description (row 19, col C) - add ", unknown" before the
#    closing parenthesis/ellipsis.
# ---------------------------------------------------------------------------
$wsAttr.Range("C19").Value = "Information about the isolation source (i.e. blood, laboratory experiment, urine, unknown...)"

# ---------------------------------------------------------------------------
# 2. collection_date description (row 24, col C) - rewrite with rich text:
#    bold the literal format tokens (YYYY-MM-DD, YYYY-MM, YYYY, unknown) and
#    add the new "or unknown" clause.
# ---------------------------------------------------------------------------
$dateCell = $wsAttr.Range("C24")
$dateText = "The date of the sample collection. Use one of the following format: YYYY-MM-DD, YYYY-MM or YYYY or unknown."
$dateCell.Value = $dateText

$dateCell.Characters(69, 10).Font.Bold = $true   # YYYY-MM-DD
$dateCell.Characters(79, 2).Font.Bold = $false   # ", "
$dateCell.Characters(81, 7).Font.Bold = $true    # YYYY-MM
$dateCell.Characters(88, 4).Font.Bold = $false   # " or "
$dateCell.Characters(92, 4).Font.Bold = $true    # YYYY
$dateCell.Characters(96, 4).Font.Bold = $false   # " or "
$dateCell.Characters(100, 7).Font.Bold = $true   # unknown
$dateCell.Characters(107, 1).Font.Bold = $false  # "."

# ---------------------------------------------------------------------------
# 3. New column D next to the collection_date row: widen it and give D24 a
#    bold, wrap-text style (matching the new font/cell style added upstream).
# ---------------------------------------------------------------------------
$wsAttr.Columns.Item(4).ColumnWidth = 28.4
$wsAttr.Range("D24").Font.Bold = $true
$wsAttr.Range("D24").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Refresh the view: scroll the "Attribute description" sheet down and move
#    the selection to C26 (last row), then re-select B7 on "Metadata".
# ---------------------------------------------------------------------------
$wsAttr.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsAttr.Range("C26").Select()

# ---------------------------------------------------------------------------
# 5. Save the workbook on completion (per commit message).
# ---------------------------------------------------------------------------
$wb.Save()
